# "success change data and store in .xlsx"
#
# RealTimePlayerData sheet: the player picked a devil boss, so the game's
# save/export step rewrote the runtime data table:
#   - Speed (row 4, col A) changed from 2 to 0.5
#   - a new ChooseDevil column (H) was appended, carrying which boss the
#     player chose (BoneMan), with its English/Chinese header and its
#     declared "string" type, matching the existing header/type rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Player speed stat changed.
$ws.Range("A4").Value = 0.5

# New "ChooseDevil" column: English header / Chinese header / declared
# type / chosen value, mirroring the layout of the existing columns.
$ws.Range("H1").Value = "ChooseDevil"
$ws.Range("H2").Value = "玩家選擇的魔王"
$ws.Range("H3").Value = "string"
$ws.Range("H4").Value = "BoneMan"

# Match the formatting already used by the other header/type/value cells
# (row 2's style) instead of inventing a new style entry.
$ws.Range("A2").Copy()
$ws.Range("H1:H4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Touching the page setup forces the worksheet to carry an explicit
# headerFooter element, matching the re-saved file's structure.
$ws.PageSetup.CenterHeader = $ws.PageSetup.CenterHeader
